# Slide 12, shape 2 ("Marcador de Posição de Conteúdo 2") - last run of the
# single paragraph is rewritten: "Estruturas Hierárquicas" and "loops" become
# italic runs, and "factorização" is corrected to "factoração".

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(12)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# The paragraph currently ends with a single run holding everything from
# ", porém ..." through "... factorização.". Replace just that run's text
# (keeps it as one run, no stray splits) with the corrected wording.
$run5 = $tr.Runs($tr.Runs().Count)
$run5.Text = ", porém há alguns casos onde a recursão pode tornar o código mais simples e claro, por exemplo em Estruturas Hierárquicas onde usar a recursão ao invés de loops torna o código mais legível e intuitivo, ou também em problemas matemáticos recursivos como por exemplo a factoração."

# Now italicise "Estruturas Hierárquicas" and "loops" using absolute
# character offsets into the whole text range (this splits the run into
# five: plain / italic / plain / italic / plain).
$tr.Characters(245, 23).Font.Italic = -1
$tr.Characters(302, 5).Font.Italic = -1
